# Update odds data on Sheet1 to reflect refreshed FlashScore figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Union de Santa Fe vs Racing Club)
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8

# Row 5 (Penarol vs Liverpool M.)
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("T5").Value = 1.8
$ws.Range("W5").Value = 3.5
$ws.Range("X5").Value = 1.29
$ws.Range("AA5").Value = 2
$ws.Range("AB5").Value = 1.73
$ws.Range("AJ5").Value = 7.5
$ws.Range("AR5").Value = 41
